# Updates cryptocurrency price & 1h-volume-change figures in the table (row 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to remain plain text so values such as "605.49" are not
# silently reinterpreted as numbers by Excel (matches source data which stores them as text).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.864.47"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "3.324.67"
$ws.Range("E3").Value = "  +2.52%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "605.49"
$ws.Range("E5").Value = "  +2.11%  "
$ws.Range("D6").Value = "143.03"
$ws.Range("E6").Value = "  +1.03%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.323.74"
$ws.Range("E8").Value = "  +2.74%  "
$ws.Range("D9").Value = "0.520"
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").Value = "0.150"
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("D11").Value = "5.57"
$ws.Range("E11").Value = "  +4.29%  "
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("D13").Value = "0.0000248"
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("D14").Value = "35.11"
$ws.Range("E14").Value = "  +2.16%  "
$ws.Range("D15").Value = "3.873.82"
$ws.Range("E15").Value = "  +2.56%  "
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").Value = "3.328.52"
$ws.Range("E17").Value = "  +2.50%  "
$ws.Range("D18").Value = "63.988.31"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").Value = "6.87"
$ws.Range("E19").Value = "  +1.37%  "
$ws.Range("D20").Value = "482.70"
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("D21").Value = "14.10"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "0.740"
$ws.Range("E22").Value = "  +2.47%  "
$ws.Range("D23").Value = "7.99"
$ws.Range("E23").Value = "  +1.02%  "
$ws.Range("D24").Value = "13.99"
$ws.Range("E24").Value = "  +6.08%  "
$ws.Range("D25").Value = "85.00"
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("E27").Value = "  +2.01%  "
$ws.Range("D28").Value = "8.30"
$ws.Range("E28").Value = "  +3.00%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "7.20"
$ws.Range("E30").Value = "  -3.22%  "
$ws.Range("D31").Value = "2.16"
$ws.Range("E31").Value = "  +2.60%  "
$ws.Range("D32").Value = "28.91"
$ws.Range("E32").Value = "  +5.06%  "
$ws.Range("E33").Value = "  -1.58%  "
$ws.Range("D34").Value = "2.53"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").Value = "1.10"
$ws.Range("E35").Value = "  +1.25%  "
$ws.Range("D36").Value = "6.09"
$ws.Range("E36").Value = "  +3.34%  "
$ws.Range("D37").Value = "0.0₃0752"
$ws.Range("E37").Value = "  +5.67%  "
$ws.Range("D38").Value = "52.39"
$ws.Range("E38").Value = "  -0.59%  "
$ws.Range("E39").Value = "  +1.89%  "
$ws.Range("D40").Value = "435.49"
$ws.Range("E40").Value = "  +3.38%  "
$ws.Range("D41").Value = "3.132.82"
$ws.Range("E41").Value = "  +5.17%  "
$ws.Range("D42").Value = "0.118"
$ws.Range("E42").Value = "  +7.21%  "
$ws.Range("E43").Value = "  +0.52%  "
$ws.Range("D44").Value = "8.37"
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("D45").Value = "0.268"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").Value = "2.27"
$ws.Range("E46").Value = "  +4.82%  "
$ws.Range("D47").Value = "36.92"
$ws.Range("E47").Value = "  +9.35%  "
$ws.Range("D48").Value = "26.47"
$ws.Range("E48").Value = "  +2.39%  "
$ws.Range("D50").Value = "2.31"
$ws.Range("E50").Value = "  -0.72%  "
$ws.Range("D51").Value = "124.67"
$ws.Range("E51").Value = "  +2.80%  "

# Restore the default (unstyled) cell style now that the text values are safely stored,
# so the Price column keeps looking like the rest of the sheet.
$ws.Range("D2:D51").Style = "Normal"
